# Re-imported order data: phone numbers in column D are rewritten as
# right-aligned text (leading "8" instead of "7") instead of raw numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phones = @{
    2  = "89805485327"
    3  = "89009326524"
    4  = "89157462939"
    5  = "89616171966"
    6  = "89611841424"
    7  = "89003075527"
    8  = "89525559858"
    9  = "89507656669"
    10 = "89802436816"
    11 = "89155823772"
    12 = "89081304789"
    13 = "89204619277"
    14 = "89065809980"
    15 = "89515672194"
    16 = "89204422077"
    17 = "89511549281"
}

$noAlignRows = @(14, 17)

foreach ($row in 2..17) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    if ($noAlignRows -notcontains $row) {
        $cell.HorizontalAlignment = -4152
    }
    $cell.Value = $phones[$row]
}

[void]$ws.Range("D17").Select()
